$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the three face-movement sheets (OralOcular_Right / OcularOral_Left
#    / OcularOral_Right shuffle). Order matters to avoid name collisions, so
#    use temporary names first, then apply the final names.
# ---------------------------------------------------------------------------
$wsOralOcularRight = $wb.Worksheets.Item("OralOcular_Right")
$wsOcularOralLeft  = $wb.Worksheets.Item("OcularOral_Left")
$wsOcularOralRight = $wb.Worksheets.Item("OcularOral_Right")

$wsOralOcularRight.Name = "__tmp1__"
$wsOcularOralLeft.Name  = "__tmp2__"
$wsOcularOralRight.Name = "__tmp3__"

$wsOralOcularRight.Name = "OralOcular_Left"
$wsOcularOralLeft.Name  = "OralOcular_Right"
$wsOcularOralRight.Name = "OcularOral_Left"

# ---------------------------------------------------------------------------
# 2. "OralOcular_Left" sheet (formerly "OralOcular_Right"):
#    row2 IMG_0490/Yes/No -> IMG_9330/Yes/No ; row3 (IMG_5694) removed.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OralOcular_Left")
$ws.Range("A2").Value = "IMG_9330"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Rows(3).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. "OralOcular_Right" sheet (formerly "OcularOral_Left"):
#    row2 IMG_3812/Yes/No -> IMG_4923/Yes/No
#    row3 IMG_4210/No/Yes -> IMG_5694/Yes/No
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OralOcular_Right")
$ws.Range("A2").Value = "IMG_4923"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Range("A3").Value = "IMG_5694"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "No"

# ---------------------------------------------------------------------------
# 4. "OcularOral_Left" sheet (formerly "OcularOral_Right"):
#    row2 IMG_7365/Yes/No -> IMG_2068/Yes/No
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OcularOral_Left")
$ws.Range("A2").Value = "IMG_2068"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"

# ---------------------------------------------------------------------------
# 5. "SnarlSmile_Left":
#    row2 IMG_0495/No/Yes  -> IMG_0504/Yes/No
#    row3 IMG_2814/Yes/No  -> IMG_7365/Yes/No
#    row4 IMG_5694/No/Yes  -> IMG_8514/Yes/No
#    row5 (new)             IMG_9374/Yes/No
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SnarlSmile_Left")
$ws.Range("A2").Value = "IMG_0504"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Range("A3").Value = "IMG_7365"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "No"
$ws.Range("A4").Value = "IMG_8514"
$ws.Range("B4").Value = "Yes"
$ws.Range("C4").Value = "No"
$ws.Range("A5").Value = "IMG_9374"
$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = "No"

# ---------------------------------------------------------------------------
# 6. "SnarlSmile_Right":
#    row2 IMG_0495/Yes/No -> IMG_4157/Yes/No
#    rows 3-6 (IMG_2122, IMG_2501, IMG_4210, IMG_7365) removed
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SnarlSmile_Right")
$ws.Range("A2").Value = "IMG_4157"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Rows("3:6").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 7. "Mentalis_Left":
#    row2 IMG_2259/No/Yes -> IMG_0490/Yes/No
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mentalis_Left")
$ws.Range("A2").Value = "IMG_0490"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"

# ---------------------------------------------------------------------------
# 8. "Mentalis_Right":
#    row2 IMG_1339/Yes/No unchanged
#    row3 IMG_2737/Yes/No -> IMG_2814/Yes/No
#    row4 IMG_3102/Yes/No -> IMG_3324/No/Yes
#    row5 IMG_3324/No/Yes -> IMG_8537/Yes/No
#    rows 6-9 (IMG_3812, IMG_4923, IMG_5694, IMG_7365) removed
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mentalis_Right")
$ws.Range("A3").Value = "IMG_2814"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "No"
$ws.Range("A4").Value = "IMG_3324"
$ws.Range("B4").Value = "No"
$ws.Range("C4").Value = "Yes"
$ws.Range("A5").Value = "IMG_8537"
$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = "No"
$ws.Rows("6:9").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 9. "Hypertonicity_Left":
#    row2 IMG_3170/Yes/No -> IMG_2737/Yes/No
#    row3 IMG_5694/Yes/No unchanged
#    row4 (new) IMG_9640/No/Yes
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hypertonicity_Left")
$ws.Range("A2").Value = "IMG_2737"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Range("A4").Value = "IMG_9640"
$ws.Range("B4").Value = "No"
$ws.Range("C4").Value = "Yes"

# ---------------------------------------------------------------------------
# 10. "Hypertonicity_Right":
#     row2 IMG_3847/Yes/No -> IMG_1339/Yes/No
#     row3 IMG_8537/Yes/No -> IMG_3847/Yes/No
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Hypertonicity_Right")
$ws.Range("A2").Value = "IMG_1339"
$ws.Range("B2").Value = "Yes"
$ws.Range("C2").Value = "No"
$ws.Range("A3").Value = "IMG_3847"
$ws.Range("B3").Value = "Yes"
$ws.Range("C3").Value = "No"
